# Creacion de documento Maestro de muestreo de datos, arreglo de muestreos en
# centro de evento y pagos, y arreglo de diagramas.
#
# This script rebuilds the "Oferta" worksheet: inserts a short-name column
# (B) and a new "Evento" column (G), fixes the trailing concatenation
# formula + hyperlinks to track the shift, and tweaks a couple of selection
# / formatting leftovers on the other worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Oferta sheet: insert the two new columns.
# ---------------------------------------------------------------------
$oferta = $wb.Worksheets.Item("Oferta")

# Insert new column B ("Nombre" corto de la oferta) - pushes old B..F to C..G
$oferta.Columns.Item(2).Insert()
# Insert new column G ("Evento") - pushes old F (now G) to H
$oferta.Columns.Item(7).Insert()

# --- New column B content -------------------------------------------------
$oferta.Range("B1").Value = "Nombre"
$oferta.Range("B2").Value = "OfertaPedicuraMadres"
$oferta.Range("B3").Value = "OfertaFacial"
$oferta.Range("B4").Value = "OfertaMasaje"

# --- New column G content --------------------------------------------------
$oferta.Range("G1").Value = "Evento"
$oferta.Range("G2").Value = "OfertaPedicura Madres"
$oferta.Range("G3").Value = "OfertaFacial San Valentin"
$oferta.Range("G4").Value = "OfertaMasaje Black Friday"

# Column insert leaves the new G cells carrying the date-format style that
# used to live in old column F; restore plain bordered formatting (matching
# column B / A) by pasting format from A2:A4.
$oferta.Range("A2:A4").Copy() | Out-Null
$oferta.Range("G2:G4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Fix the trailing "display" formula (now column H) --------------------
# Was: =B#&"-"&C#   (old Servicio name & old Descuento)
# Now: =B#&" "&D#   (new short Nombre & Descuento, now column D)
$oferta.Range("H2").Formula = "=B2&"" ""&D2"
$oferta.Range("H3").Formula = "=B3&"" ""&D3"
$oferta.Range("H4").Formula = "=B4&"" ""&D4"

# --- Fix hyperlinks: they used to anchor on column B, which is now C ------
$oferta.Hyperlinks.Delete()
$oferta.Hyperlinks.Add($oferta.Range("C2"), "", "TipoServicio!B2", [System.Reflection.Missing]::Value, "TipoServicio!B2") | Out-Null
$oferta.Hyperlinks.Add($oferta.Range("C3:C4"), "", "TipoServicio!B2", [System.Reflection.Missing]::Value, "TipoServicio!B2") | Out-Null

# --- Column widths for the two new columns ---------------------------------
$oferta.Columns.Item(2).ColumnWidth = 20.1666667   # -> stored width 21
$oferta.Columns.Item(7).ColumnWidth = 23.4523809   # -> stored width ~24.285

# --- Stray underline-styled empty cell left by the author at H10 ----------
$oferta.Range("H10").Font.Underline = 2
$oferta.Range("H10").Value = ""

# --- Selection left on the sheet -------------------------------------------
$oferta.Range("H10").Select()

# ---------------------------------------------------------------------
# 2. TipoServicio / Servicio sheets: selection cleared to the used range.
# ---------------------------------------------------------------------
$tipoServicio = $wb.Worksheets.Item("TipoServicio")
$tipoServicio.Range("A1:C4").Select()

$servicio = $wb.Worksheets.Item("Servicio")
$servicio.Range("A1:E4").Select()

# Leave the "Oferta" sheet active/selected, matching the saved workbook view.
$oferta.Activate()
$oferta.Range("H10").Select()
